$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells P1/Q1: copy the formatting (bold, border, centered) from
# the existing last header cell (O1) so the new columns match the header
# row's style, then set their values.
$ws.Range("O1").Copy() | Out-Null
$ws.Range("P1:Q1").PasteSpecial(-4122) | Out-Null
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15

# For each data row (2..25):
#  - swap I/K values (I: 1 -> 2, K: 2 -> 1)
#  - swap M/O values (M: 1 -> 2, O: 2 -> 1)
#  - add new P, Q columns with value 2 (no special style, like the other
#    un-styled data columns B..O)
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 9).Value = 2   # I
    $ws.Cells.Item($r, 11).Value = 1  # K
    $ws.Cells.Item($r, 13).Value = 2  # M
    $ws.Cells.Item($r, 15).Value = 1  # O
    $ws.Cells.Item($r, 16).Value = 2  # P
    $ws.Cells.Item($r, 17).Value = 2  # Q
}
